$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F20").Value = 1097
$wsExhibit.Range("F25").Value = 1384
$wsExhibit.Range("F26").Value = 4029
$wsExhibit.Range("F28").Value = 43521
$wsExhibit.Range("F32").Value = 642
$wsExhibit.Range("F41").Value = 835
$wsExhibit.Range("F46").Value = 63

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 224

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F8").Value = 2295
$wsLocal.Range("F9").Value = 9176
$wsLocal.Range("F10").Value = 1441

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1441
$wsAll.Range("F17").Value = 224
$wsAll.Range("F21").Value = 1097
$wsAll.Range("F25").Value = 4029
$wsAll.Range("F39").Value = 835
$wsAll.Range("F43").Value = 63
